$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2 through 362) holds the "Förändrad" (last changed) date.
# Update every value from serial 45203 (2023-10-04) to serial 45204 (2023-10-05).
for ($r = 2; $r -le 362; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
